{"js": "const pairs = [\n  [\"2025-10-03 Friday\", \"2025-10-04 Saturday\"],\n  [\"188\u00f75=\", \"548\u00f75=\"],\n  [\"909\u00f76=\", \"693\u00f77=\"],\n  [\"827\u00f73=\", \"216\u00f79=\"],\n  [\"983\u00f72=\", \"248\u00f77=\"],\n  [\"612\u00f72=\", \"510\u00f78=\"],\n  [\"622\u00f74=\", \"477\u00f76=\"],\n  [\"365\u00f75=\", \"385\u00f73=\"],\n  [\"341\u00f74=\", \"279\u00f72=\"],\n  [\"667\u00f76=\", \"702\u00f74=\"],\n  [\"443\u00f76=\", \"921\u00f74=\"],\n  [\"180\u00f79=\", \"352\u00f78=\"],\n  [\"163\u00f76=\", \"739\u00f78=\"],\n  [\"662\u00f74=\", \"599\u00f75=\"],\n  [\"356\u00f72=\", \"786\u00f78=\"],\n  [\"662\u00f77=\", \"767\u00f76=\"],\n  [\"161\u00f79=\", \"440\u00f74=\"],\n  [\"552\u00f75=\", \"195\u00f79=\"],\n  [\"505\u00f72=\", \"533\u00f74=\"],\n  [\"424\u00f75=\", \"247\u00f78=\"],\n  [\"405\u00f74=\", \"546\u00f73=\"],\n  [\"792\u00f79=\", \"718\u00f75=\"],\n  [\"921\u00f76=\", \"801\u00f77=\"],\n  [\"978\u00f73=\", \"986\u00f79=\"],\n  [\"924\u00f78=\", \"318\u00f75=\"],\n  [\"538\u00f77=\", \"302\u00f72=\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-10-03 Friday\", \"2025-10-04 Saturday\"),\n    @(\"188\u00f75=\", \"548\u00f75=\"),\n    @(\"909\u00f76=\", \"693\u00f77=\"),\n    @(\"827\u00f73=\", \"216\u00f79=\"),\n    @(\"983\u00f72=\", \"248\u00f77=\"),\n    @(\"612\u00f72=\", \"510\u00f78=\"),\n    @(\"622\u00f74=\", \"477\u00f76=\"),\n    @(\"365\u00f75=\", \"385\u00f73=\"),\n    @(\"341\u00f74=\", \"279\u00f72=\"),\n    @(\"667\u00f76=\", \"702\u00f74=\"),\n    @(\"443\u00f76=\", \"921\u00f74=\"),\n    @(\"180\u00f79=\", \"352\u00f78=\"),\n    @(\"163\u00f76=\", \"739\u00f78=\"),\n    @(\"662\u00f74=\", \"599\u00f75=\"),\n    @(\"356\u00f72=\", \"786\u00f78=\"),\n    @(\"662\u00f77=\", \"767\u00f76=\"),\n    @(\"161\u00f79=\", \"440\u00f74=\"),\n    @(\"552\u00f75=\", \"195\u00f79=\"),\n    @(\"505\u00f72=\", \"533\u00f74=\"),\n    @(\"424\u00f75=\", \"247\u00f78=\"),\n    @(\"405\u00f74=\", \"546\u00f73=\"),\n    @(\"792\u00f79=\", \"718\u00f75=\"),\n    @(\"921\u00f76=\", \"801\u00f77=\"),\n    @(\"978\u00f73=\", \"986\u00f79=\"),\n    @(\"924\u00f78=\", \"318\u00f75=\"),\n    @(\"538\u00f77=\", \"302\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
